# "Some more styling and HTML"
#
# For each of the three card sheets (Blue/sheet1, Yellow/sheet2, Green/sheet3):
#   1. Insert a new row 1 with totals: G1 = SUM(G3:G<last>), and H1:K1 a
#      shared formula SUM(H3:H<last>) filled right across H..K.
#   2. Rename the constructor used in the generated "const XXnn = new
#      ....Card(...)" line (column M) from the colour-specific class name
#      (BlueCard / YellowCard / GreenCard) to the shared "Card" name.
#   3. Rewrite the generated closing-row formula (column B, last row) from
#      `return [...]` to `activeDeck.splice(0, oldCards, ...)`.
#   4. (Green sheet only) drop the now-unused helper column L (row sum of
#      G:K) that isn't used anywhere else.
#
# Finally re-select the "current cell" on each sheet (B<last>) so the
# selection matches, doing Green last so it remains the active tab (as it
# was originally).

$wb = $excel.ActiveWorkbook

$sheets = @(
    @{ Index = 1; LastDataRow = 22; LastRow = 24; OldClass = "BlueCard"   },
    @{ Index = 2; LastDataRow = 32; LastRow = 34; OldClass = "YellowCard" },
    @{ Index = 3; LastDataRow = 42; LastRow = 44; OldClass = "GreenCard"  }
)

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Index)
    $lastData = $s.LastDataRow
    $lastRow  = $s.LastRow

    # --- 1. New header/total row -------------------------------------------------
    $ws.Range("G1").Formula    = "=SUM(G3:G$lastData)"
    $ws.Range("H1:K1").Formula = "=SUM(H3:H$lastData)"

    # --- 2. Constructor rename in column M --------------------------------------
    $ws.Range("M3").Formula = '=CONCAT("const ",C3," = new Card(",D3,$A$3,E3,$A$3,F3,$A$3,G3,$A$3,H3,$A$3,I3,$A$3,J3,$A$3,K3,")")'
    $mRange = "M4:M$lastData"
    $ws.Range($mRange).FormulaR1C1 = '=CONCAT("const ",RC[-10]," = new Card(",RC[-9],R3C1,RC[-8],R3C1,RC[-7],R3C1,RC[-6],R3C1,RC[-5],R3C1,RC[-4],R3C1,RC[-3],R3C1,RC[-2],")")'

    # --- 3. Sheet-specific clean-up (Green: drop helper column L) --------------
    if ($s.OldClass -eq "GreenCard") {
        $ws.Range("L3:L$lastData").ClearContents()
    }

    # --- 4. Closing-row formula ---------------------------------------------------
    $bRange = "B3:B$lastData"
    $ws.Range("B$lastRow").Formula = '=CONCAT("activeDeck.splice(0, oldCards, ",' + $bRange + ',")")'
}

# --- 5. Selection: B<last> on every sheet, Green selected/activated last so it
#        remains the active tab, matching the workbook's original state.
foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Index)
    [void]$ws.Range("B$($s.LastRow)").Select()
}
